$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new value would otherwise be
# auto-detected as a number by Excel, so they stay text like the source data.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated values
$ws.Range("D2").Value = "29.359.64"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").Value = "1.876.36"
$ws.Range("E3").Value = "  +0.46%  "
$ws.Range("D4").Value = "0.9998"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "0.7129"
$ws.Range("E5").Value = "  -1.87%  "
$ws.Range("D6").Value = "241.92"
$ws.Range("E6").Value = "  +0.35%  "
$ws.Range("D8").Value = "0.3115"
$ws.Range("E8").Value = "  +0.66%  "
$ws.Range("D9").Value = "0.07715"
$ws.Range("E9").Value = "  -2.18%  "
$ws.Range("D10").Value = "25.16"
$ws.Range("E10").Value = "  -0.40%  "
$ws.Range("D11").Value = "0.08376"
$ws.Range("E11").Value = "  +1.56%  "
$ws.Range("D12").Value = "1.890.71"
$ws.Range("E12").Value = "  +1.72%  "
$ws.Range("D13").Value = "5.235"
$ws.Range("E13").Value = "  -0.16%  "
$ws.Range("D14").Value = "0.7140"
$ws.Range("E14").Value = "  -1.31%  "
$ws.Range("D15").Value = "91.79"
$ws.Range("E15").Value = "  +1.10%  "
$ws.Range("D16").Value = "29.367.40"
$ws.Range("E16").Value = "  +0.05%  "
$ws.Range("D17").Value = "0.000008301"
$ws.Range("E17").Value = "  +5.94%  "
$ws.Range("D18").Value = "5.961"
$ws.Range("E18").Value = "  +1.65%  "
$ws.Range("D19").Value = "243.23"
$ws.Range("E19").Value = "  -0.32%  "
$ws.Range("D20").Value = "2.135.48"
$ws.Range("E20").Value = "  +0.53%  "
$ws.Range("D21").Value = "13.23"
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("D22").Value = "0.9992"
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("D23").Value = "7.891"
$ws.Range("E23").Value = "  -1.61%  "
$ws.Range("D24").Value = "1.000"
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("D25").Value = "0.1631"
$ws.Range("E25").Value = "  +1.90%  "
$ws.Range("D26").Value = "163.75"
$ws.Range("E26").Value = "  +0.89%  "
$ws.Range("D27").Value = "9.028"
$ws.Range("E27").Value = "  +0.75%  "
$ws.Range("D28").Value = "18.56"
$ws.Range("E29").Value = "  +0.59%  "
$ws.Range("D30").Value = "4.409"
$ws.Range("E30").Value = "  +0.12%  "
$ws.Range("D31").Value = "1.288"
$ws.Range("E31").Value = "  -4.59%  "
$ws.Range("E32").Value = "  +4.99%  "
$ws.Range("D33").Value = "0.05243"
$ws.Range("E33").Value = "  +0.55%  "
$ws.Range("D34").Value = "1.931"
$ws.Range("E34").Value = "  -0.29%  "
$ws.Range("D35").Value = "0.7554"
$ws.Range("E35").Value = "  +3.64%  "
$ws.Range("D36").Value = "1.177"
$ws.Range("E36").Value = "  -0.84%  "
$ws.Range("D37").Value = "2.681"
$ws.Range("E37").Value = "  +0.13%  "
$ws.Range("D38").Value = "0.01862"
$ws.Range("E38").Value = "  +0.24%  "
$ws.Range("D39").Value = "2.725"
$ws.Range("E39").Value = "  +0.83%  "
$ws.Range("D40").Value = "1.159.24"
$ws.Range("E40").Value = "  -0.47%  "
$ws.Range("D41").Value = "6.361"
$ws.Range("E41").Value = "  +4.19%  "
$ws.Range("D42").Value = "73.33"
$ws.Range("E42").Value = "  +1.37%  "
$ws.Range("D43").Value = "0.8890"
$ws.Range("E43").Value = "  -1.82%  "
$ws.Range("D44").Value = "104.80"
$ws.Range("D45").Value = "0.9994"
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("D46").Value = "2.031.83"
$ws.Range("E46").Value = "  +0.65%  "
$ws.Range("D47").Value = "1.798"
$ws.Range("E47").Value = "  +1.17%  "
$ws.Range("E48").Value = "  -1.55%  "
$ws.Range("D49").Value = "9.419"
$ws.Range("E49").Value = "  +1.52%  "
$ws.Range("D50").Value = "0.4308"
$ws.Range("E50").Value = "  +0.76%  "
$ws.Range("D51").Value = "7.042"
$ws.Range("E51").Value = "  +0.16%  "

# Restore default cell style (clears the transient text-format flag)
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
